$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "2012-06-02"
}
$rng.ClearFormats()
